$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two hex operand cells used by the calculator (row 9).
# Write C9 first then B9 so the shared-string table order matches
# the target ("d000" gets the lower index than "d800").
$ws.Range("C9").Value = "'d000"
$ws.Range("B9").Value = "d800"

# Un-hide the helper rows (10-17); rows 11-13 have no content, so
# delete+insert them to drop their row entry entirely instead of
# leaving a bare <row/> behind.
$ws.Rows("10:17").Hidden = $false
$ws.Rows(11).EntireRow.Delete()
$ws.Rows(11).EntireRow.Insert()
$ws.Rows(12).EntireRow.Delete()
$ws.Rows(12).EntireRow.Insert()
$ws.Rows(13).EntireRow.Delete()
$ws.Rows(13).EntireRow.Insert()

# Move the active selection to C9 (reflecting where the user was working).
$ws.Range("C9").Select()
